$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("B2", 0.8122020464430761),
    @("C2", 0.134595722647191),
    @("D2", 0.1849994499776031),
    @("E2", 0.1509367944984383),
    @("F2", 1.307652727885412),
    @("J2", 0.1621527824637177),
    @("M2", 0.3355715691277226),
    @("N2", 1.262359827795187),
    @("O2", 3.087754504164849),
    @("B3", 0.7328332203227035),
    @("C3", 0.117459001797414),
    @("D3", 0.1833800931512855),
    @("E3", 0.1506820903074768),
    @("F3", 1.305795992869363),
    @("J3", 0.1626666169639464),
    @("M3", 0.3150397459795542),
    @("N3", 1.275168404928138),
    @("O3", 3.087176000199094),
    @("B4", 0.6842055932467019),
    @("C4", 0.1069039566617391),
    @("D4", 0.1824459978109445),
    @("E4", 0.1505882780729699),
    @("F4", 1.305437182570287),
    @("J4", 0.1630602117257034),
    @("M4", 0.3025399794181638),
    @("N4", 1.283539258625105),
    @("O4", 3.088817609988411),
    @("B5", 0.6644169043132706),
    @("C5", 0.102594617345602),
    @("D5", 0.1820805468703028),
    @("E5", 0.150565808663611),
    @("F5", 1.305487348835115),
    @("J5", 0.1632402463020846),
    @("M5", 0.2974733748442588),
    @("N5", 1.287077808195257),
    @("O5", 3.089988244028802),
    @("B6", 0.6611326965133912),
    @("C6", 0.1018785741293016),
    @("D6", 0.1820207839667702),
    @("E6", 0.1505630302076781),
    @("F6", 1.305507538229925),
    @("J6", 0.1632713272994089),
    @("M6", 0.2966337172371283),
    @("N6", 1.287673074963372),
    @("O6", 3.090212913113248),
    @("B7", 0.6839386035113932),
    @("C7", 0.1068458717814167),
    @("D7", 0.1824410075810476),
    @("E7", 0.1505879111991995),
    @("F7", 1.305437064062453),
    @("J7", 0.1630625602027571),
    @("M7", 0.3024715391060084),
    @("N7", 1.28358646502258),
    @("O7", 3.088831366950131),
    @("B8", 0.7848144384463467),
    @("C8", 0.1286939708747354),
    @("D8", 0.1844286393235279),
    @("E8", 0.150835998751937),
    @("F8", 1.306850338438394),
    @("J8", 0.1623137448550551),
    @("M8", 0.3284701807485249),
    @("N8", 1.266671156984074),
    @("O8", 3.087140316565126),
    @("B9", 0.9834310707374243),
    @("C9", 0.1712677543411019),
    @("D9", 0.1888015522033726),
    @("E9", 0.151818134162788),
    @("F9", 1.315825453781059),
    @("J9", 0.1614649933582868),
    @("M9", 0.3802916761950002),
    @("N9", 1.237516110772781),
    @("O9", 3.099692961937137),
    @("B10", 1.12980988716015),
    @("C10", 0.2023739132535241),
    @("D10", 0.1923012054911339),
    @("E10", 0.1528409420400223),
    @("F10", 1.326211335182563),
    @("J10", 0.1612193798457255),
    @("M10", 0.4188671138120057),
    @("N10", 1.218541556089214),
    @("O10", 3.118630717820793),
    @("B11", 1.196494336769376),
    @("C11", 0.2164858720063023),
    @("D11", 0.1939550189973858),
    @("E11", 0.1533714777816968),
    @("F11", 1.331761850221326),
    @("J11", 0.1611897731024428),
    @("M11", 0.4365235089979933),
    @("N11", 1.210440145693198),
    @("O11", 3.129364915390624),
    @("B12", 1.221758975671207),
    @("C12", 0.2218239880407111),
    @("D12", 0.1945901053940702),
    @("E12", 0.1535817417944401),
    @("F12", 1.333982585385627),
    @("J12", 0.1611903729858355),
    @("M12", 0.4432248629628148),
    @("N12", 1.207448586409193),
    @("O12", 3.133735064077626),
    @("B13", 1.216317237252724),
    @("C13", 0.2206745902231262),
    @("D13", 0.1944529367853818),
    @("E13", 0.1535360416652907),
    @("F13", 1.333499021585141),
    @("J13", 0.1611897184468702),
    @("M13", 0.4417809326093618),
    @("N13", 1.208089479688297),
    @("O13", 3.132780284921125),
    @("B14", 1.19857262598714),
    @("C14", 0.2169251592593469),
    @("D14", 0.1940070914761662),
    @("E14", 0.1533885888331454),
    @("F14", 1.331942168375605),
    @("J14", 0.1611895857278327),
    @("M14", 0.4370745293209239),
    @("N14", 1.210192499603842),
    @("O14", 3.129718327312446),
    @("B15", 1.187705158037375),
    @("C15", 0.2146277646717181),
    @("D15", 0.1937351454785272),
    @("E15", 0.1532994882260823),
    @("F15", 1.331004035156624),
    @("J15", 0.1611910426747443),
    @("M15", 0.4341936997467286),
    @("N15", 1.211490593335377),
    @("O15", 3.127882571649337),
    @("B16", 1.12545373749407),
    @("C16", 0.2014508691613912),
    @("D16", 0.1921943630473493),
    @("E16", 0.1528075815039287),
    @("F16", 1.325865224984582),
    @("J16", 0.1612229674294383),
    @("M16", 0.4177153788507084),
    @("N16", 1.219081673563124),
    @("O16", 3.11797190894174),
    @("B17", 1.087288296636984),
    @("C17", 0.1933572557215371),
    @("D17", 0.1912649266208888),
    @("E17", 0.1525225104456815),
    @("F17", 1.322924348012634),
    @("J17", 0.1612635875585369),
    @("M17", 0.4076339680218837),
    @("N17", 1.223874371549449),
    @("O17", 3.11243526347792),
    @("B18", 1.065345646410435),
    @("C18", 0.1886984252591049),
    @("D18", 0.190736160454918),
    @("E18", 0.1523646885631749),
    @("F18", 1.32131056690146),
    @("J18", 0.1612946809672522),
    @("M18", 0.4018456128867243),
    @("N18", 1.226680908694902),
    @("O18", 3.109450176644003),
    @("B19", 1.057917839449715),
    @("C19", 0.1871204151584323),
    @("D19", 0.1905581310496842),
    @("E19", 0.1523123085999636),
    @("F19", 1.320777515849585),
    @("J19", 0.1613065362003567),
    @("M19", 0.3998875374042967),
    @("N19", 1.227639724314038),
    @("O19", 3.108473714587063),
    @("B20", 1.091350139063707),
    @("C20", 0.1942192089349817),
    @("D20", 0.1913632647204793),
    @("E20", 0.152552221082221),
    @("F20", 1.323229363468144),
    @("J20", 0.1612584634590419),
    @("M20", 0.408706097362554),
    @("N20", 1.223359015614705),
    @("O20", 3.11300400275681),
    @("B21", 1.203784313084952),
    @("C21", 0.2180266168302296),
    @("D21", 0.1941378081048413),
    @("E21", 0.1534316454070215),
    @("F21", 1.332396226789328),
    @("J21", 0.1611893041385954),
    @("M21", 0.4384565023099327),
    @("N21", 1.209572721532176),
    @("O21", 3.130609406952829),
    @("B22", 1.277339966275463),
    @("C22", 0.2335523460451725),
    @("D22", 0.1960025307822519),
    @("E22", 0.1540609566875588),
    @("F22", 1.339080249918709),
    @("J22", 0.1612129508689932),
    @("M22", 0.4579889198584866),
    @("N22", 1.201007150546609),
    @("O22", 3.143895588780737),
    @("B23", 1.238075580807845),
    @("C23", 0.225269148244962),
    @("D23", 0.1950026114233765),
    @("E23", 0.153720096727568),
    @("F23", 1.335449420411166),
    @("J23", 0.161194030048641),
    @("M23", 0.4475560753049805),
    @("N23", 1.205538069296843),
    @("O23", 3.136641437087292),
    @("B24", 1.089513782602921),
    @("C24", 0.1938295376376971),
    @("D24", 0.191318788683347),
    @("E24", 0.1525387699969123),
    @("F24", 1.323091226235562),
    @("J24", 0.1612607559538759),
    @("M24", 0.4082213640474635),
    @("N24", 1.223591848697843),
    @("O24", 3.112746258992985),
    @("B25", 0.9296177607383243),
    @("C25", 0.1597801857510888),
    @("D25", 0.1875679833430439),
    @("E25", 0.1514994746219926),
    @("F25", 1.312732273432545),
    @("J25", 0.1616282430514531),
    @("M25", 0.366183764255851),
    @("N25", 1.244973677795059),
    @("O25", 3.094594007563018)
)

foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value = [double]$pair[1]
}

Write-Host "Updated $($updates.Count) cells"
